$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.980.51"
$ws.Range("E2").Value = "  +1.97%  "
$ws.Range("D3").Value = "3.708.20"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "614.46"
$ws.Range("E5").Value = "  +6.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.14"
$ws.Range("E6").Value = "  +9.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.634"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.712"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.08"
$ws.Range("E11").Value = "  +7.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000288"
$ws.Range("E12").Value = "  -3.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.50"
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("D14").Value = "4.285.34"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "3.704.78"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.127"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.23"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.13"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.87"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("D20").Value = "68.749.91"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "409.45"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.60"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "88.90"
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.03"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.73"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.87"
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.04"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("E28").Value = "  -2.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.61"
$ws.Range("E29").Value = "  +1.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.95"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.21"
$ws.Range("E31").Value = "  -9.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.58"
$ws.Range("E32").Value = "  -0.90%  "
$ws.Range("E33").Value = "  +3.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "623.56"
$ws.Range("E34").Value = "  +3.96%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "44.29"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.57"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.410"
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").Value = "0.0₃0807"
$ws.Range("E39").Value = "  -11.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.140"
$ws.Range("E41").Value = "  +3.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.01"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0441"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.60"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("E45").Value = "  +3.24%  "
$ws.Range("D46").Value = "2.860.42"
$ws.Range("E46").Value = "  +4.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.72"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.04"
$ws.Range("E48").Value = "  -4.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.07"
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.45"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.76"
$ws.Range("E51").Value = "  +0.06%  "
